$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted before the existing row 388,
# pushing every subsequent row (old 388..500) down by one (new 389..501).
$ws.Rows("388:388").Insert()

$ws.Range("A388").Value = 4
$ws.Range("B388").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C388").Value = "Los Lagos"
$ws.Range("D388").Value = 45093
$ws.Range("E388").Value = 10
$ws.Range("F388").Value = 100114014
$ws.Range("G388").Value = "Betarraga"
$ws.Range("H388").Value = "Sin especificar"
$ws.Range("I388").Value = "Primera"
$ws.Range("J388").Value = 1200
$ws.Range("K388").Value = 1000
$ws.Range("L388").Value = 1000
$ws.Range("M388").Value = 1000
$ws.Range("N388").Value = "$/paquete 5 unidades"
$ws.Range("O388").Value = "Región Metropolitana"
$ws.Range("P388").Value = 200
$ws.Range("Q388").Value = 5
$ws.Range("R388").Value = "Hortaliza"
